$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose data (columns B..AC) need to be swapped with each other.
# Column A (the running index) stays put on each row; only the match
# details in B:AC move between the two rows of each pair.
# (Using two parallel arrays instead of an array-of-arrays because the
# interpreter's `foreach` flattens nested arrays.)
$rows1 = @(7, 11, 13, 15, 30, 65, 70, 74, 78, 84, 95, 103, 106)
$rows2 = @(8, 12, 14, 16, 31, 66, 71, 75, 79, 85, 96, 104, 107)

for ($i = 0; $i -lt $rows1.Length; $i++) {
    $r1 = $rows1[$i]
    $r2 = $rows2[$i]

    $range1 = $ws.Range("B${r1}:AC${r1}")
    $range2 = $ws.Range("B${r2}:AC${r2}")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
